# Moves the "MCT-3A-Microcontroladores" class entries from column B (segunda)
# to column D (quarta) for rows 4, 6 and 7, clears it from B3 entirely,
# and adds it to D8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$placeholder = "[-, -, 'MCT-3A-Microcontroladores', -]"

# Row 3: remove the class from B3 (no replacement elsewhere in this row)
$ws.Range("B3").Value = "-"

# Row 4: move class from B4 to D4
$ws.Range("B4").Value = "-"
$ws.Range("D4").Value = $placeholder

# Row 6: move class from B6 to D6
$ws.Range("B6").Value = "-"
$ws.Range("D6").Value = $placeholder

# Row 7: move class from B7 to D7
$ws.Range("B7").Value = "-"
$ws.Range("D7").Value = $placeholder

# Row 8: add the class to D8
$ws.Range("D8").Value = $placeholder
